# Add a new "NC_PR" column to both worksheets, inserted right before the
# existing "ExF_MA" column (old column L). That pushes ExF_MA -> M and
# MA_PR -> N on both "det_full" and "det_short" sheets, and populates the
# new column with the NC_PR series.

$wb = $excel.ActiveWorkbook

# New NC_PR values for "det_full" (sheet1), keyed by row number (row 2 = year 2015 ... row 31 = year 2044)
$fullNCPR = @{
    2  = 28.233351474359893
    3  = 26.834264792199534
    4  = 25.93219265135643
    5  = 25.077051653222373
    6  = 24.25923049405211
    7  = 23.785146651455086
    8  = 23.68823758525744
    9  = 23.60056645324508
    10 = 23.51953724194511
    11 = 23.44217035343518
    12 = 23.3671070308337
    13 = 23.312198344243047
    14 = 23.263882615521844
    15 = 23.218591402403977
    16 = 23.172250822721104
    17 = 23.12171489691118
    18 = 23.083072685749684
    19 = 23.063830279312587
    20 = 23.066955719892295
    21 = 23.094765407588657
    22 = 23.148986663035576
    23 = 23.190462502906957
    24 = 23.21770097078383
    25 = 23.23199203759563
    26 = 23.232307808184533
    27 = 23.226882343745693
    28 = 23.222376915918186
    29 = 23.219641831854
    30 = 23.217573735946473
    31 = 23.21963912026572
}

# New NC_PR values for "det_short" (sheet2), keyed by row number
# (years 2015, 2020, 2025, 2030, 2035, 2040, 2044)
$shortNCPR = @{
    2 = 28.233351474359893
    3 = 23.785146651455086
    4 = 23.3671070308337
    5 = 23.12171489691118
    6 = 23.148986663035576
    7 = 23.226882343745693
    8 = 23.21963912026572
}

$sheets = @(
    @{ Name = "det_full";  Values = $fullNCPR;  LastRow = 31 },
    @{ Name = "det_short"; Values = $shortNCPR; LastRow = 8 }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Insert a new blank column at L, shifting old L (ExF_MA) -> M and old M (MA_PR) -> N
    $ws.Columns("L:L").Insert()

    # Header for the new column
    $ws.Range("L1").Value = "NC_PR"

    # Populate the new column's data values
    $values = $info.Values
    foreach ($row in $values.Keys) {
        $ws.Cells.Item($row, 12).Value = $values[$row]
    }
}
